{"js": "// Update the Java stack-trace text embedded in the document body to reflect\n// line-number shifts (POI 3.17.0 -> 4.0.1) and the swap of the Maven/Tycho/\n// Equinox launcher tail frames for the Eclipse JDT JUnit runner frames.\n\n// Small line-number replacements (each string is unique in the document).\nconst replacements = [\n  [\n    \"JavaMethodService.internalInvoke(JavaMethodService.java:163)\",\n    \"JavaMethodService.internalInvoke(JavaMethodService.java:162)\"\n  ],\n  [\n    \"AbstractService.invoke(AbstractService.java:136)\",\n    \"AbstractService.invoke(AbstractService.java:135)\"\n  ],\n  [\n    \"EvaluationServices.call(EvaluationServices.java:168)\",\n    \"EvaluationServices.call(EvaluationServices.java:172)\"\n  ],\n  [\n    \"EvaluationServices.callOrApply(EvaluationServices.java:204)\",\n    \"EvaluationServices.callOrApply(EvaluationServices.java:208)\"\n  ],\n  [\n    \"AstSwitch.doSwitch(AstSwitch.java:118)\",\n    \"AstSwitch.doSwitch(AstSwitch.java:119)\"\n  ],\n  [\n    \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)\",\n    \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:462)\"\n  ],\n  [\n    \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)\",\n    \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:372)\"\n  ],\n  [\n    \"sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)\",\n    \"sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '\" + oldText + \"' but found \" + results.items.length);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Replace the whole Maven/Tycho/Equinox launcher tail of the stack trace\n// with the Eclipse JDT JUnit runner tail.\nconst oldTail =\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)\\n\" +\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)\\n\" +\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)\\n\" +\n  \"\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:156)\\n\" +\n  \"\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)\\n\" +\n  \"\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)\\n\" +\n  \"\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)\\n\" +\n  \"\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)\\n\" +\n  \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)\\n\" +\n  \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)\\n\" +\n  \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)\\n\" +\n  \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)\\n\";\n\nconst newTail =\n  \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\\n\";\n\nconst tailResults = body.search(oldTail, { matchCase: true });\ntailResults.load(\"items\");\nawait context.sync();\nif (tailResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the Maven/Tycho/Equinox tail but found \" + tailResults.items.length);\n}\ntailResults.items[0].insertText(newTail, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the Java stack-trace text embedded in the document body to reflect\n# line-number shifts (POI 3.17.0 -> 4.0.1) and the swap of the Maven/Tycho/\n# Equinox launcher tail frames for the Eclipse JDT JUnit runner frames.\n\n$d = $word.ActiveDocument\n\n# Small line-number replacements (each string is unique in the document).\n$replacements = @(\n  @(\"JavaMethodService.internalInvoke(JavaMethodService.java:163)\", \"JavaMethodService.internalInvoke(JavaMethodService.java:162)\"),\n  @(\"AbstractService.invoke(AbstractService.java:136)\", \"AbstractService.invoke(AbstractService.java:135)\"),\n  @(\"EvaluationServices.call(EvaluationServices.java:168)\", \"EvaluationServices.call(EvaluationServices.java:172)\"),\n  @(\"EvaluationServices.callOrApply(EvaluationServices.java:204)\", \"EvaluationServices.callOrApply(EvaluationServices.java:208)\"),\n  @(\"AstSwitch.doSwitch(AstSwitch.java:118)\", \"AstSwitch.doSwitch(AstSwitch.java:119)\"),\n  @(\"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)\", \"AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:462)\"),\n  @(\"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)\", \"AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:372)\"),\n  @(\"sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)\", \"sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\")\n)\n\nforeach ($pair in $replacements) {\n  $range = $d.Content\n  $range.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n\n# Replace the whole Maven/Tycho/Equinox launcher tail of the stack trace\n# with the Eclipse JDT JUnit runner tail.\n$oldTail = \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n\" +\n  \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n\" +\n  \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n\" +\n  \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" +\n  \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" +\n  \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" +\n  \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" +\n  \"`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n\" +\n  \"`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n\" +\n  \"`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n\" +\n  \"`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n\" +\n  \"`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n\" +\n  \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" +\n  \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" +\n  \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" +\n  \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" +\n  \"`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n\" +\n  \"`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n\" +\n  \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n\" +\n  \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n\" +\n  \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n\" +\n  \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n\" +\n  \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" +\n  \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" +\n  \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" +\n  \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" +\n  \"`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n\" +\n  \"`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n\" +\n  \"`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n\" +\n  \"`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)`n\"\n\n$newTail = \"`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n\" +\n  \"`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n\" +\n  \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n\" +\n  \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n\" +\n  \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n\" +\n  \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)`n\"\n\n$range = $d.Content\n$range.Find.Execute($oldTail, $false, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)\n"}
